$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# Remove the duplicated "Contact" row (old row 11); everything below shifts up by one.
$meta.Rows.Item(11).Delete()

# Update Version
$meta.Range("B3").Value = "6.0.0"

# Update Date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now has a value
$meta.Range("B9").Value = "Alvearie Team"

# Row 10 changes from Contact/No-display-for-ContactDetail to Jurisdiction/United States of America
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Elements sheet ---
$elem = $wb.Worksheets.Item("Elements")

# Row 2 (the top-level "Extension" row) gets a specific Short/Definition instead of the generic boilerplate
$elem.Range("K2").Value = "Long Term Disability Benefit Rate"
$elem.Range("L2").Value = "Long term disability (LTD) benefit rate, expressed as a percentage of base wages (for example, 50.00%, 60.00%, 70.00%)"
